$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 10021
$ws.Range("I20").Value = 10021
$ws.Range("K20").Value = 10021
$ws.Range("M20").Value = -9791
$ws.Range("H35").Value = 10021
$ws.Range("I35").Value = 10021
$ws.Range("K35").Value = 10021
$ws.Range("M35").Value = -9642
$ws.Range("H76").Value = 3283.3333
$ws.Range("I76").Value = 3282.353
$ws.Range("K76").Value = 3282.353
$ws.Range("M76").Value = -2967.353
$ws.Range("H79").Value = 3283.3333
$ws.Range("I79").Value = 3282.353
$ws.Range("K79").Value = 3282.353
$ws.Range("M79").Value = -2190.353
$ws.Range("H88").Value = 8358.714
$ws.Range("I88").Value = 8501.5
$ws.Range("J88").Value = 8301.6
$ws.Range("K88").Value = 8501.5
$ws.Range("L88").Value = 8301.6
$ws.Range("M88").Value = -8095.5
$ws.Range("N88").Value = -9113.6
$ws.Range("H91").Value = 8358.714
$ws.Range("I91").Value = 8501.5
$ws.Range("J91").Value = 8301.6
$ws.Range("K91").Value = 8501.5
$ws.Range("L91").Value = 8301.6
$ws.Range("M91").Value = -7097.5
$ws.Range("N91").Value = -11109.6
$ws.Range("H132").Value = 20204106
$ws.Range("I132").Value = 2584985.5
$ws.Range("J132").Value = 83339290
$ws.Range("K132").Value = 7754956.5
$ws.Range("L132").Value = 250017870
$ws.Range("M132").Value = -7752426.5
$ws.Range("N132").Value = -250022930
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5695.6
$ws.Range("I32").Value = 3298.8982
$ws.Range("K32").Value = 3298.8982
$ws.Range("M32").Value = -3011.8982
$ws.Range("H88").Value = 2450.3333
$ws.Range("I88").Value = 2493.3333
$ws.Range("J88").Value = 2407.3333
$ws.Range("K88").Value = 2493.3333
$ws.Range("L88").Value = 2407.3333
$ws.Range("M88").Value = -2087.3333
$ws.Range("N88").Value = -3219.3333
$ws.Range("H91").Value = 2450.3333
$ws.Range("I91").Value = 2493.3333
$ws.Range("J91").Value = 2407.3333
$ws.Range("K91").Value = 2493.3333
$ws.Range("L91").Value = 2407.3333
$ws.Range("M91").Value = -1089.3333
$ws.Range("N91").Value = -5215.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1018.55
$ws.Range("I94").Value = 935.0625
$ws.Range("J94").Value = 1352.5
$ws.Range("K94").Value = 935.0625
$ws.Range("L94").Value = 1352.5
$ws.Range("M94").Value = -484.0625
$ws.Range("N94").Value = -2254.5
$ws.Range("H99").Value = 1157.9
$ws.Range("I99").Value = 1157
$ws.Range("J99").Value = 1160
$ws.Range("K99").Value = 1157
$ws.Range("L99").Value = 1160
$ws.Range("M99").Value = 341
$ws.Range("N99").Value = -4156
$ws.Range("H134").Value = 4239.6665
$ws.Range("I134").Value = 3611.6667
$ws.Range("J134").Value = 4993.2666
$ws.Range("K134").Value = 10835.0001
$ws.Range("L134").Value = 14979.7998
$ws.Range("M134").Value = -8300.000100000001
$ws.Range("N134").Value = -20049.7998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 998.875
$ws.Range("I16").Value = 740
$ws.Range("J16").Value = 1085.1666
$ws.Range("K16").Value = 740
$ws.Range("L16").Value = 1085.1666
$ws.Range("M16").Value = -453
$ws.Range("N16").Value = -1659.1666
$ws.Range("H113").Value = 998.875
$ws.Range("I113").Value = 740
$ws.Range("J113").Value = 1085.1666
$ws.Range("K113").Value = 740
$ws.Range("L113").Value = 1085.1666
$ws.Range("M113").Value = 1430
$ws.Range("N113").Value = -5425.1666
$ws.Range("H132").Value = 2873.5518
$ws.Range("I132").Value = 2364.2856
$ws.Range("J132").Value = 4210.375
$ws.Range("K132").Value = 7092.8568
$ws.Range("L132").Value = 12631.125
$ws.Range("M132").Value = -4562.8568
$ws.Range("N132").Value = -17691.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 860.27
$ws.Range("I131").Value = 568.9
$ws.Range("J131").Value = 892.6445
$ws.Range("K131").Value = 1706.7
$ws.Range("L131").Value = 2677.9335
$ws.Range("M131").Value = 3333.3
$ws.Range("N131").Value = -12757.9335
$ws.Range("H132").Value = 4136519.8
$ws.Range("I132").Value = 1756298.9
$ws.Range("J132").Value = 22226198
$ws.Range("K132").Value = 15806690.1
$ws.Range("L132").Value = 200035782
$ws.Range("M132").Value = -15804160.1
$ws.Range("N132").Value = -200040842
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3335.1292
$ws.Range("I132").Value = 3178.1853
$ws.Range("J132").Value = 4394.5
$ws.Range("K132").Value = 9534.555899999999
$ws.Range("L132").Value = 13183.5
$ws.Range("M132").Value = -7004.555899999999
$ws.Range("N132").Value = -18243.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 934.2632
$ws.Range("I46").Value = 1268.5
$ws.Range("J46").Value = 780
$ws.Range("K46").Value = 1268.5
$ws.Range("L46").Value = 780
$ws.Range("M46").Value = -1080.5
$ws.Range("N46").Value = -1156
$ws.Range("H61").Value = 1230.9166
$ws.Range("I61").Value = 1103.9615
$ws.Range("J61").Value = 1561
$ws.Range("K61").Value = 1103.9615
$ws.Range("L61").Value = 1561
$ws.Range("M61").Value = -901.9614999999999
$ws.Range("N61").Value = -1965
$ws.Range("H113").Value = 1230.9166
$ws.Range("I113").Value = 1103.9615
$ws.Range("J113").Value = 1561
$ws.Range("K113").Value = 1103.9615
$ws.Range("L113").Value = 1561
$ws.Range("M113").Value = 1066.0385
$ws.Range("N113").Value = -5901
$ws.Range("H132").Value = 3156.8604
$ws.Range("I132").Value = 2914.6287
$ws.Range("J132").Value = 4216.625
$ws.Range("K132").Value = 8743.8861
$ws.Range("L132").Value = 12649.875
$ws.Range("M132").Value = -6213.8861
$ws.Range("N132").Value = -17709.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 31000
$ws.Range("J46").Value = 31000
$ws.Range("L46").Value = 31000
$ws.Range("N46").Value = -31462
$ws.Range("H92").Value = 29108.166
$ws.Range("J92").Value = 29108.166
$ws.Range("L92").Value = 29108.166
$ws.Range("N92").Value = -34100.166
$ws.Range("H128").Value = 35454.547
$ws.Range("J128").Value = 35454.547
$ws.Range("L128").Value = 35454.547
$ws.Range("N128").Value = -45414.547
$ws.Range("H132").Value = 4765187.5
$ws.Range("I132").Value = 6669996
$ws.Range("J132").Value = 3166.9167
$ws.Range("K132").Value = 20009988
$ws.Range("L132").Value = 9500.750100000001
$ws.Range("M132").Value = -20007458
$ws.Range("N132").Value = -14560.7501
$ws.Range("H134").Value = 31000
$ws.Range("J134").Value = 31000
$ws.Range("L134").Value = 93000
$ws.Range("N134").Value = -98070
